$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28, shifting existing rows 28:57 down to 29:58
$ws.Rows("28:28").Insert()

# Fill in the new weekly record at row 28
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = [DateTime]::FromOADate(44789)
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112035
$ws.Range("G28").Value = "Bruselas (repollito)"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = 21000
$ws.Range("L28").Value = 21000
$ws.Range("M28").Value = 21000
$ws.Range("N28").Value = "`$/malla 15 kilos"
$ws.Range("O28").Value = "Hijuelas"
$ws.Range("P28").Value = 1400
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = "Hortaliza"
